# Added player set 2
# Appends two new result rows (rows 6 & 7) to Sheet1, extending the
# used range from A1:C5 to A1:C7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Tue_Dec__5_12_34_52_2023"
$ws.Range("B6").Value = "f"
$ws.Range("C6").Value = 30

$ws.Range("A7").Value = "Tue_Dec__5_12_37_31_2023"
$ws.Range("B7").Value = "f"
$ws.Range("C7").Value = 30
